$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52 - this pushes the existing rows 52:120
# down to 53:121 (matches the new dimension A1:R121 and the shifted data
# seen across the rest of the sheet).
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with this week's record (same
# market/category/quality metadata as its neighbours, new date + prices).
$ws.Range("A52").Value = 2
$ws.Range("B52").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C52").Value = "Coquimbo"
$ws.Range("D52").Value = 45225
$ws.Range("E52").Value = 4
$ws.Range("F52").Value = 100112026
$ws.Range("G52").Value = "Haba"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 900
$ws.Range("K52").Value = 7000
$ws.Range("L52").Value = 8000
$ws.Range("M52").Value = 7500
$ws.Range("N52").Value = "$/saco 25 kilos"
$ws.Range("O52").Value = "Provincia de Limarí"
$ws.Range("P52").Value = 300
$ws.Range("Q52").Value = 25
$ws.Range("R52").Value = "Hortaliza"
